# Observer pattern: append new "booked" order rows (14-17) to the orders sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(14, "booked", 75.97999999999999, 1),
    @(15, "booked", 64.99, 1),
    @(16, "booked", 44.99, 1),
    @(17, "booked", 0.0, 1)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
